# Delete column I (duplicate of column H) and shift J:L left to I:K,
# matching the "다중 꺾은 선 그래프" data-prep cleanup described in the
# commit message: the stray duplicate column that used to sit at I
# (an exact copy of H) is removed, and the downstream columns collapse
# left by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("I").Delete()

# Widen column G a touch (13 characters), and leave the selection where
# Excel naturally lands after a column delete: on the column that slid
# into the deleted one's place.
$ws.Columns("G").ColumnWidth = 86/7
$ws.Range("I1:I1048576").Select() | Out-Null
